# Wordless acknowledgments sheet: remove the "lxml" dependency row.
#
# The "lxml" library entry occupied row 11 (Name/Home Page/Version/Authors/
# License/License URL). Deleting it shifts every later row up by one and
# also invalidates the worksheet's hyperlink collection (each hyperlink is
# anchored to an absolute cell), so we rebuild the hyperlinks for the new
# layout afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the lxml row entirely - everything below shifts up by one.
$ws.Rows(11).Delete()

# 2) The hyperlink collection does not auto-repair on row delete, so clear
#    it out and re-create the links against their new (post-shift) cells,
#    in the same relative order as before (minus the two links that
#    belonged to the deleted lxml row).
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.crummy.com/software/BeautifulSoup/")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/Ousret/charset_normalizer")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/Mimino666/langdetect")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://github.com/saffsd/langid.py")
$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.numpy.org/")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://foss.heptapod.net/openpyxl/openpyxl")
$ws.Hyperlinks.Add($ws.Range("B25"), "https://github.com/python-openxml/python-docx")
$ws.Hyperlinks.Add($ws.Range("B27"), "https://github.com/psf/requests")
$ws.Hyperlinks.Add($ws.Range("B29"), "https://scipy.org/scipylib/")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://bazaar.launchpad.net/~leonardr/beautifulsoup/bs4/view/head:/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/Ousret/charset_normalizer/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://github.com/Mimino666/langdetect/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://github.com/saffsd/langid.py/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://github.com/numpy/numpy/blob/master/LICENSE.txt")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://foss.heptapod.net/openpyxl/openpyxl/-/blob/branch/3.0/LICENCE.rst")
$ws.Hyperlinks.Add($ws.Range("F25"), "https://github.com/python-openxml/python-docx/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F27"), "https://github.com/requests/requests/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F29"), "https://github.com/scipy/scipy/blob/master/LICENSE.txt")
$ws.Hyperlinks.Add($ws.Range("F24"), "https://docs.python.org/3.8/license.html", "psf-license-agreement-for-python-release")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://github.com/pyinstaller/pyinstaller/blob/develop/COPYING.txt")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://www.riverbankcomputing.com/static/Docs/PyQt5/introduction.html", "license")
$ws.Hyperlinks.Add($ws.Range("B24"), "https://www.python.org/")
$ws.Hyperlinks.Add($ws.Range("B19"), "http://www.pyinstaller.org/")
$ws.Hyperlinks.Add($ws.Range("B22"), "https://riverbankcomputing.com/software/pyqt/")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://matplotlib.org/")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://networkx.org/")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://matplotlib.org/users/license.html")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://github.com/networkx/networkx/blob/master/LICENSE.txt")
$ws.Hyperlinks.Add($ws.Range("B36"), "https://github.com/amueller/word_cloud")
$ws.Hyperlinks.Add($ws.Range("F36"), "https://github.com/amueller/word_cloud/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/Esukhia/botok")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/cltk/cltk")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/fxsjy/jieba")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://github.com/taishi-i/nagisa")
$ws.Hyperlinks.Add($ws.Range("B14"), "http://www.nltk.org/")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://github.com/yichen0831/opencc-python")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://github.com/lancopku/pkuseg-python")
$ws.Hyperlinks.Add($ws.Range("B21"), "https://pyphen.org/")
$ws.Hyperlinks.Add($ws.Range("B20"), "https://github.com/kmike/pymorphy2")
$ws.Hyperlinks.Add($ws.Range("B23"), "https://github.com/PyThaiNLP/pythainlp")
$ws.Hyperlinks.Add($ws.Range("B26"), "https://github.com/natasha/razdel")
$ws.Hyperlinks.Add($ws.Range("B28"), "https://github.com/alvations/sacremoses")
$ws.Hyperlinks.Add($ws.Range("B30"), "https://spacy.io/")
$ws.Hyperlinks.Add($ws.Range("B31"), "https://github.com/ponrawee/ssg")
$ws.Hyperlinks.Add($ws.Range("B33"), "https://github.com/sloria/TextBlob")
$ws.Hyperlinks.Add($ws.Range("B34"), "https://github.com/mideind/Tokenizer")
$ws.Hyperlinks.Add($ws.Range("B35"), "https://github.com/undertheseanlp/underthesea")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/Esukhia/botok/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/cltk/cltk/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://github.com/fxsjy/jieba/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://github.com/taishi-i/nagisa/blob/master/LICENSE.txt")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://github.com/nltk/nltk/blob/develop/LICENSE.txt")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://github.com/yichen0831/opencc-python/blob/master/LICENSE.txt")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://github.com/lancopku/pkuseg-python/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://github.com/Kozea/Pyphen/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://github.com/kmike/pymorphy2/", "pymorphy2")
$ws.Hyperlinks.Add($ws.Range("F23"), "https://github.com/PyThaiNLP/pythainlp/blob/dev/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F26"), "https://github.com/natasha/razdel/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F28"), "https://github.com/alvations/sacremoses/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F30"), "https://github.com/explosion/spaCy/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F31"), "https://github.com/ponrawee/ssg/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F33"), "https://github.com/sloria/TextBlob/blob/dev/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F34"), "https://github.com/mideind/Tokenizer/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F35"), "https://github.com/undertheseanlp/underthesea/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/Xangis/extra-stopwords")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://github.com/michmech/lemmatization-lists")
$ws.Hyperlinks.Add($ws.Range("B32"), "https://github.com/stopwords-iso/stopwords-iso")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://github.com/Xangis/extra-stopwords/blob/master/LICENSE")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://github.com/michmech/lemmatization-lists/blob/master/LICENCE")
$ws.Hyperlinks.Add($ws.Range("F32"), "https://github.com/stopwords-iso/stopwords-iso/blob/master/LICENSE")

# 3) Restore the view state: frozen pane back at the top of the data and
#    the active selection on the now-empty-ish area where row 11 used to be.
$ws.Range("B2").Select()
$ws.Range("A9").Select()
